# Edit script implementing the "meta-release" casing fix described in the
# commit message: "always 'meta-release' except in heading or at start of
# sentence."
#
# Slide-index -> content mapping (confirmed via SlideID):
#   Slide 2 (SlideID 2044) -> "Release Management / Motivation and Benefits"
#   Slide 3 (SlideID 2047) -> "Release Management / Release Cycle"
#   Slide 5 (SlideID 2081) -> "Beyond Fall24 Meta-release - Upcoming APIs"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - "Motivation and Benefits": shape id 7 "Inhaltsplatzhalter 6"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shContent = $s2.Shapes.Item(3)

# Explicitly pin the placeholder's position/size (was inheriting from the
# layout; author dragged/confirmed it to an explicit frame).
$shContent.Left = 428045 / 12700.0
$shContent.Top = 1890346 / 12700.0
$shContent.Width = 10944225 / 12700.0
$shContent.Height = 4124383 / 12700.0

$tr2 = $shContent.TextFrame.TextRange
$full2 = $tr2.Text
$pos = 0
while ($true) {
    $found = $full2.IndexOf("CAMARA Meta-release", $pos)
    if ($found -eq -1) { break }
    $isPlural = $full2.Substring($found, 21) -eq "CAMARA Meta-releases "
    if ($isPlural) {
        $tr2.Characters($found + 1, 21).Text = "CAMARA meta-releases "
    } else {
        $tr2.Characters($found + 1, 20).Text = "CAMARA meta-release "
    }
    $pos = $found + 1
}

# ---------------------------------------------------------------------
# Slide 3 - "Release Cycle"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# "Fall / Meta- / release" and "Spring / Meta- / release" ovals, nested
# inside the "Group 34" group shape.
$grp = $s3.Shapes.Item(4)
$ovalFall = $grp.GroupItems.Item(3)
$trFall = $ovalFall.TextFrame.TextRange
$fallFull = $trFall.Text
$fallIdx = $fallFull.IndexOf("Meta-")
$trFall.Characters($fallIdx + 1, 5).Text = "meta-"

$ovalSpring = $grp.GroupItems.Item(4)
$trSpring = $ovalSpring.TextFrame.TextRange
$springFull = $trSpring.Text
$springIdx = $springFull.IndexOf("Meta-")
$trSpring.Characters($springIdx + 1, 5).Text = "meta-"

# "TextBox 60" (id 61) - "...of CAMARA Fall Meta-Release"
$tbFall = $s3.Shapes.Item(13)
$tbFall.Left = 9124406 / 12700.0
$tbFall.Top = 5077821 / 12700.0
$tbFall.Width = 2254143 / 12700.0
$tbFall.Height = 769441 / 12700.0
$trTbFall = $tbFall.TextFrame.TextRange
$tbFallFull = $trTbFall.Text
$tbFallIdx = $tbFallFull.IndexOf("Fall Meta-Release")
$trTbFall.Characters($tbFallIdx + 1, 18).Text = "Fall meta-release"

# "TextBox 61" (id 62) - "...of CAMARA Spring Meta-Release"
$tbSpring = $s3.Shapes.Item(14)
$tbSpring.Left = 3851279 / 12700.0
$tbSpring.Top = 1408489 / 12700.0
$tbSpring.Width = 2473755 / 12700.0
$tbSpring.Height = 769441 / 12700.0
$trTbSpring = $tbSpring.TextFrame.TextRange
$tbSpringFull = $trTbSpring.Text
$tbSpringIdx = $tbSpringFull.IndexOf("Spring Meta-Release")
$trTbSpring.Characters($tbSpringIdx + 1, 20).Text = "Spring meta-release"

# ---------------------------------------------------------------------
# Slide 5 - "Beyond Fall24 Meta-release - Upcoming APIs" title
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$titleShape = $s5.Shapes.Item(2)
$trTitle = $titleShape.TextFrame.TextRange
$titleFull = $trTitle.Text
$breakIdx = $titleFull.IndexOf([char]11)
$trTitle.Characters(1, $breakIdx).Text = "Beyond Fall24 Meta-release "
$afterFull = $trTitle.Text
$breakIdx2 = $afterFull.IndexOf([char]11)
$startRun2 = $breakIdx2 + 2
$len2 = $afterFull.Length - $breakIdx2 - 1
$trTitle.Characters($startRun2, $len2).Text = [string]([char]8211) + " Upcoming APIs"
